$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column "groupe_id" header at J1, pushing the old J1 header
# ("gestionnaires_additionnels") to the new column K1.
$ws.Range("K1").Value = "gestionnaires_additionnels"
$ws.Range("J1").Value = "groupe_id"

# Update the selection to match the target workbook state.
$ws.Range("J1").Select()
